$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 41 (new data row, formatted like row 38) ---
$ws.Range("A38:I38").Copy()
$ws.Range("A41:I41").PasteSpecial(-4122)
$ws.Range("F28").Copy()
$ws.Range("F41").PasteSpecial(-4122)

$ws.Range("A41").Value = "5/5/2022"
$ws.Range("B41").Value = "RASY"
$ws.Range("C41").Value = "N/A"
$ws.Range("D41").Value = "Clermont"
$ws.Range("E41").Value = "Capitale-Nationale"
$ws.Range("F41").Value = "C/D"
$ws.Range("G41").Value = "Cote 3"
$ws.Range("I41").Value = "Eve Murray"

# --- Row 42 (new data row, formatted like row 39) ---
$ws.Range("A39:I39").Copy()
$ws.Range("A42:I42").PasteSpecial(-4122)
$ws.Range("F28").Copy()
$ws.Range("F42").PasteSpecial(-4122)

$ws.Range("A42").Value = "5/5/2022"
$ws.Range("B42").Value = "PSCR"
$ws.Range("C42").Value = "N/A"
$ws.Range("D42").Value = "Clermont"
$ws.Range("E42").Value = "Capitale-Nationale"
$ws.Range("F42").Value = "C/D"
$ws.Range("G42").Value = "Cote 1"
$ws.Range("I42").Value = "Eve Murray"

$excel.CutCopyMode = $false

# Matches the final selection recorded in the workbook after the edit.
$ws.Range("C47").Select()
